# Update the "想去人数" (want-to-go count) figures in column F for a handful
# of rows on both the "展览" and "全部类型" worksheets. These two sheets carry
# duplicate data, so the same edits must be applied to both.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new F-column value
$updates = @{
    2  = 1766
    7  = 12183
    11 = 425
    14 = 13562
    23 = 2115
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
